$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.171343088150024
$ws.Range("B1").Value = 2.593072891235352
$ws.Range("C1").Value = 9.455362319946289
$ws.Range("D1").Value = 2.100996971130371
$ws.Range("E1").Value = 1.228919267654419
